# Automatische test-sync: 2025-08-04 20:31:50
$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# Append the new mail log entry as row 15
$row = 15
$logs.Cells.Item($row, 1).Value = "Kun jij dit afhandelen?"
$logs.Cells.Item($row, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($row, 3).Value = "Testmail #3: Kun jij dit afhandelen?"
$logs.Cells.Item($row, 4).Value = "Planning / Afspraak"
$logs.Cells.Item($row, 5).Value = "Bedankt, we hebben dit doorgestuurd naar planning@bedrijf.nl."
$logs.Cells.Item($row, 6).Value = "2025-08-04 20:31:25"
$logs.Cells.Item($row, 7).Value = "Ja"
$logs.Cells.Item($row, 8).Value = "Ja"
$logs.Cells.Item($row, 9).Value = "Nee"
$logs.Cells.Item($row, 10).Value = "Nee"

# Update the Dashboard count for "Planning / Afspraak" (5 -> 6)
$dashboard.Range("B2").Value = 6

# Expand conditional formatting ranges to include the newly added row 15
$cfColumns = @("D", "G", "H", "I", "J")
foreach ($col in $cfColumns) {
    $oldRange = $logs.Range("$($col)2:$($col)14")
    $newRange = $logs.Range("$($col)2:$($col)15")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}
